$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').NumberFormat = '@'
$ws.Range('D2').Value = '58.078.53'
$ws.Range('D3').NumberFormat = '@'
$ws.Range('D3').Value = '2.357.52'
$ws.Range('E3').NumberFormat = '@'
$ws.Range('E3').Value = '  +1.20%  '
$ws.Range('E4').NumberFormat = '@'
$ws.Range('E4').Value = '  -0.17%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '545.03'
$ws.Range('E5').NumberFormat = '@'
$ws.Range('E5').Value = '  +0.84%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '134.84'
$ws.Range('E6').NumberFormat = '@'
$ws.Range('E6').Value = '  -0.57%  '
$ws.Range('E7').NumberFormat = '@'
$ws.Range('E7').Value = '  +0.72%  '
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '0.558'
$ws.Range('E8').NumberFormat = '@'
$ws.Range('E8').Value = '  +4.30%  '
$ws.Range('E9').NumberFormat = '@'
$ws.Range('E9').Value = '  +0.01%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '5.59'
$ws.Range('E10').NumberFormat = '@'
$ws.Range('E10').Value = '  +3.38%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '0.151'
$ws.Range('E11').NumberFormat = '@'
$ws.Range('E11').Value = '  -1.33%  '
$ws.Range('E12').NumberFormat = '@'
$ws.Range('E12').Value = '  +0.57%  '
$ws.Range('B13').NumberFormat = '@'
$ws.Range('B13').Value = 'WrappedliquidstakedEther2.0'
$ws.Range('C13').NumberFormat = '@'
$ws.Range('C13').Value = 'https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth'
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '2.776.26'
$ws.Range('E13').NumberFormat = '@'
$ws.Range('E13').Value = '  +0.76%  '
$ws.Range('B14').NumberFormat = '@'
$ws.Range('B14').Value = 'Avalanche'
$ws.Range('C14').NumberFormat = '@'
$ws.Range('C14').Value = 'https://coinranking.com/coin/dvUj0CzDZ+avalanche-avax'
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '23.82'
$ws.Range('E14').NumberFormat = '@'
$ws.Range('E14').Value = '  +0.43%  '
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '58.023.46'
$ws.Range('E15').NumberFormat = '@'
$ws.Range('E15').Value = '  +1.05%  '
$ws.Range('E16').NumberFormat = '@'
$ws.Range('E16').Value = '  +0.63%  '
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '2.363.13'
$ws.Range('E17').NumberFormat = '@'
$ws.Range('E17').Value = '  +0.34%  '
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '10.79'
$ws.Range('E18').NumberFormat = '@'
$ws.Range('E18').Value = '  +2.80%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '331.52'
$ws.Range('E19').NumberFormat = '@'
$ws.Range('E19').Value = '  -1.76%  '
$ws.Range('E20').NumberFormat = '@'
$ws.Range('E20').Value = '  +1.90%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '6.70'
$ws.Range('E21').NumberFormat = '@'
$ws.Range('E21').Value = '  -1.76%  '
$ws.Range('E22').NumberFormat = '@'
$ws.Range('E22').Value = '  +0.20%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '62.68'
$ws.Range('E23').NumberFormat = '@'
$ws.Range('E23').Value = '  +1.84%  '
$ws.Range('E24').NumberFormat = '@'
$ws.Range('E24').Value = '  -1.25%  '
$ws.Range('B25').NumberFormat = '@'
$ws.Range('B25').Value = 'Binance-PegBSC-USD'
$ws.Range('C25').NumberFormat = '@'
$ws.Range('C25').Value = 'https://coinranking.com/coin/i5jggxiwp+binance-pegbsc-usd-bsc-usd'
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '1.00'
$ws.Range('E25').NumberFormat = '@'
$ws.Range('E25').Value = '  +1.02%  '
$ws.Range('B26').NumberFormat = '@'
$ws.Range('B26').Value = 'InternetComputer(DFINITY)'
$ws.Range('C26').NumberFormat = '@'
$ws.Range('C26').Value = 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp'
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '8.46'
$ws.Range('E26').NumberFormat = '@'
$ws.Range('E26').Value = '  -0.12%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '1.36'
$ws.Range('E27').NumberFormat = '@'
$ws.Range('E27').Value = '  -1.45%  '
$ws.Range('E28').NumberFormat = '@'
$ws.Range('E28').Value = '  +0.44%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '170.87'
$ws.Range('E29').NumberFormat = '@'
$ws.Range('E29').Value = '  -1.74%  '
$ws.Range('E30').NumberFormat = '@'
$ws.Range('E30').Value = '  +0.29%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '6.15'
$ws.Range('E31').NumberFormat = '@'
$ws.Range('E31').Value = '  +0.15%  '
$ws.Range('B32').NumberFormat = '@'
$ws.Range('B32').Value = 'SuiNetwork'
$ws.Range('C32').NumberFormat = '@'
$ws.Range('C32').Value = 'https://coinranking.com/coin/3xJluUMvp+suinetwork-sui'
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '1.03'
$ws.Range('E32').NumberFormat = '@'
$ws.Range('E32').Value = '  +2.70%  '
$ws.Range('B33').NumberFormat = '@'
$ws.Range('B33').Value = 'EthereumClassic'
$ws.Range('C33').NumberFormat = '@'
$ws.Range('C33').Value = 'https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc'
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '18.47'
$ws.Range('E33').NumberFormat = '@'
$ws.Range('E33').Value = '  -0.12%  '
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '4.26'
$ws.Range('E35').NumberFormat = '@'
$ws.Range('E35').Value = '  +4.03%  '
$ws.Range('E37').NumberFormat = '@'
$ws.Range('E37').Value = '  -2.14%  '
$ws.Range('E38').NumberFormat = '@'
$ws.Range('E38').Value = '  +1.03%  '
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '39.46'
$ws.Range('E39').NumberFormat = '@'
$ws.Range('E39').Value = '  +0.32%  '
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '143.46'
$ws.Range('E40').NumberFormat = '@'
$ws.Range('E40').Value = '  -3.91%  '
$ws.Range('B41').NumberFormat = '@'
$ws.Range('B41').Value = 'Filecoin'
$ws.Range('C41').NumberFormat = '@'
$ws.Range('C41').Value = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '3.67'
$ws.Range('E41').NumberFormat = '@'
$ws.Range('E41').Value = '  +1.18%  '
$ws.Range('B42').NumberFormat = '@'
$ws.Range('B42').Value = 'PolygonEcosystemToken'
$ws.Range('C42').NumberFormat = '@'
$ws.Range('C42').Value = 'https://coinranking.com/coin/iDZ0tG-wI+polygonecosystemtoken-pol'
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '0.380'
$ws.Range('E42').NumberFormat = '@'
$ws.Range('E42').Value = '  +0.59%  '
$ws.Range('B43').NumberFormat = '@'
$ws.Range('B43').Value = 'Bittensor'
$ws.Range('C43').NumberFormat = '@'
$ws.Range('C43').Value = 'https://coinranking.com/coin/pgv7xSFi6+bittensor-tao'
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '291.14'
$ws.Range('E43').NumberFormat = '@'
$ws.Range('E43').Value = '  +1.78%  '
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '0.0942'
$ws.Range('E44').NumberFormat = '@'
$ws.Range('E44').Value = '  +1.51%  '
$ws.Range('E45').NumberFormat = '@'
$ws.Range('E45').Value = '  +0.64%  '
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '19.06'
$ws.Range('E46').NumberFormat = '@'
$ws.Range('E46').Value = '  +0.77%  '
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '0.567'
$ws.Range('E47').NumberFormat = '@'
$ws.Range('E47').Value = '  +0.74%  '
$ws.Range('E48').NumberFormat = '@'
$ws.Range('E48').Value = '  +2.68%  '
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '17.60'
$ws.Range('E49').NumberFormat = '@'
$ws.Range('E49').Value = '  +0.29%  '
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '0.380'
$ws.Range('E50').NumberFormat = '@'
$ws.Range('E50').Value = '  -0.43%  '
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '11.08'
$ws.Range('E51').NumberFormat = '@'
$ws.Range('E51').Value = '  +1.82%  '
